# Fixed #366 User content is lost after two generation without edition.
#
# Converts the two "simple field" (<w:fldSimple>) user-content markers
# ("m:usercontent zone1" and "m:endusercontent") into equivalent
# "complex field" runs (begin fldChar / instrText / separate fldChar /
# end fldChar), each living alone in its own paragraph - matching what
# Word itself produces once such a field has been updated/regenerated.

$d = $word.ActiveDocument

function Get-ComplexFieldParagraphXml($instrText) {
    # A pkg:package payload wrapping a single <w:p> - InsertXML replaces
    # the target Range's content with this paragraph's runs.
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:instrText>' + $instrText + '</w:instrText></w:r>' +
        '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Snapshot every simple-field's instruction text and position first -
# Fields.Item(i) handles shift once earlier fields get rewritten, so we
# resolve all target paragraphs up-front before mutating anything.
$targets = @()
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $f = $d.Fields.Item($i)
    [void]($targets += ,@($f.Code.Start, $f.Code.Text.Trim()))
}

foreach ($t in $targets) {
    $fieldPos = $t[0]
    $instr = $t[1]

    # Locate the paragraph that owns this field code position.
    $target = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($fieldPos -ge $p.Range.Start -and $fieldPos -lt $p.Range.End) {
            $target = $p
            break
        }
    }

    if ($target -ne $null) {
        $xml = Get-ComplexFieldParagraphXml $instr
        [void]$target.Range.InsertXML($xml)
    }
}
